$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original inlineStr/text cells) for the numeric-looking
# values we are about to write, then strip the temporary format again so no stray
# cell style is left behind.
$fmtRange = $ws.Range("B2:G51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "235.91"
$ws.Range("G2").Value = "10"
$ws.Range("D3").Value = "22.10"
$ws.Range("G3").Value = "10"
$ws.Range("D4").Value = "5.443"
$ws.Range("G4").Value = "10"
$ws.Range("D5").Value = "0.05637"
$ws.Range("G5").Value = "10"
$ws.Range("D6").Value = "3.372"
$ws.Range("G6").Value = "10"
$ws.Range("D7").Value = "6.472"
$ws.Range("G7").Value = "10"
$ws.Range("D8").Value = "1.075"
$ws.Range("G8").Value = "10"
$ws.Range("D9").Value = "0.7879"
$ws.Range("G9").Value = "10"
$ws.Range("D10").Value = "0.1401"
$ws.Range("G10").Value = "10"
$ws.Range("D11").Value = "0.07340"
$ws.Range("G11").Value = "10"
$ws.Range("D12").Value = "0.03207"
$ws.Range("G12").Value = "10"
$ws.Range("D13").Value = "0.02976"
$ws.Range("G13").Value = "10"
$ws.Range("B14").Value = "ProBitToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D14").Value = "0.1054"
$ws.Range("E14").Value = "13ProBitTokenPROB"
$ws.Range("G14").Value = "10"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09249"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").Value = "10"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001659"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Value = "10"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "3.262"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("G17").Value = "10"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04763"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").Value = "10"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "0.0005756"
$ws.Range("E19").Value = "18OneONE"
$ws.Range("G19").Value = "10"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "0.006248"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("G20").Value = "10"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.005098"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("G21").Value = "10"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "0.001052"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("G22").Value = "10"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "0.0001501"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("G23").Value = "10"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "3.861"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("G24").Value = "10"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "2.153"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("G25").Value = "10"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "0.3291"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("G26").Value = "10"
$ws.Range("G27").Value = "10"
$ws.Range("G28").Value = "10"
$ws.Range("G29").Value = "10"
$ws.Range("G30").Value = "10"
$ws.Range("G31").Value = "10"
$ws.Range("G32").Value = "10"
$ws.Range("G33").Value = "10"
$ws.Range("G34").Value = "10"
$ws.Range("G35").Value = "10"
$ws.Range("G36").Value = "10"
$ws.Range("G37").Value = "10"
$ws.Range("G38").Value = "10"
$ws.Range("G39").Value = "10"
$ws.Range("G40").Value = "10"
$ws.Range("D41").Value = "0.006955"
$ws.Range("G41").Value = "10"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.003503"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("G42").Value = "10"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1037"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").Value = "10"
$ws.Range("D44").Value = "0.009903"
$ws.Range("G44").Value = "10"
$ws.Range("D45").Value = "0.00005446"
$ws.Range("G45").Value = "10"
$ws.Range("G46").Value = "10"
$ws.Range("D47").Value = "0.6758"
$ws.Range("G47").Value = "10"
$ws.Range("D48").Value = "0.03842"
$ws.Range("G48").Value = "10"
$ws.Range("G49").Value = "10"
$ws.Range("G50").Value = "10"
$ws.Range("G51").Value = "10"

$fmtRange.ClearFormats()
